# Apply edit: update Approved/Rejected column (I) and clear ReasonToReject
# column (J) for rows 2-9, then update the selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 9).Value = "Approved"
    $ws.Cells.Item($r, 10).ClearContents()
}

$ws.Activate()
$ws.Range("J10").Select()
